# Automatische test-sync: 2025-09-02 22:58:50
# Adds a new mail-log entry (row 3) to the "Logs" sheet and the matching
# aggregate row (row 3) to the "Dashboard" sheet, then extends the
# conditional-formatting ranges and the dashboard chart's category/value
# series so they include the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet - append the new mail entry in row 3
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Bestelling M6 bouten"
$logs.Range("B3").Value = '"Testbedrijf 123 B.V." <klantenservice@testbedrijf123.nl>'
$logs.Range("D3").Value = "Inkoop / Bestellingen"
$logs.Range("F3").Value = "2025-09-02 22:57:55"
$logs.Range("G3").Value = "Nee"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# ---------------------------------------------------------------------
# 2. "Dashboard" sheet - append the aggregated count in row 3
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Inkoop / Bestellingen"
$dash.Range("B3").Value = 1

# ---------------------------------------------------------------------
# 3. Extend conditional formatting ranges on "Logs" to include row 3
# ---------------------------------------------------------------------
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D3"))
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G3"))
$logs.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H3"))
$logs.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I3"))
$logs.Range("J2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J3"))

# ---------------------------------------------------------------------
# 4. Extend the Dashboard bar chart series (category + value) to row 3
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$3,'Dashboard'!`$B`$2:`$B`$3,1)"
